$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (oldest observation) so everything shifts up by one row,
# matching the new data vintage, then remove the now-empty last row (19).
$ws.Rows.Item(2).Delete()

# New forecast values (recomputed) for E2:E18.
$eValues = @(
    0.4944284391569687,
    -0.4782015746048418,
    1.324233212457782,
    0.7478380109886329,
    -0.2445716668737163,
    -0.2617076051026235,
    -0.100009932057743,
    0.3000376062062493,
    0.1740313431290996,
    0.3390041783450259,
    0.2210188332817387,
    -0.09571633453315798,
    -1.49562970548649,
    -0.1048501255800471,
    0.9692952624595019,
    0.1544084105021826,
    0.3997355152047577
)

for ($i = 0; $i -lt $eValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 5).Value = $eValues[$i]
}
